$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) on Price cells whose new value looks like a plain
# number, so Excel stores them as literal strings (preserving trailing zeros / exact
# formatting) instead of silently converting them to floating-point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated cell values (Coin, Link, Price, Volume(1h)).
$ws.Range("D2").Value = "68.422.00"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.289.69"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "587.74"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "186.08"
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +4.71%  "
$ws.Range("D10").Value = "6.72"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "0.420"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "3.859.48"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "28.98"
$ws.Range("E14").Value = "  +4.13%  "
$ws.Range("D15").Value = "68.408.61"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "3.288.66"
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "13.70"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "384.85"
$ws.Range("D21").Value = "7.81"
$ws.Range("D22").Value = "71.45"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +4.01%  "
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").Value = "0.195"
$ws.Range("E26").Value = "  +8.10%  "
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.32"
$ws.Range("E31").Value = "  +4.00%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "23.02"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +3.88%  "
$ws.Range("D36").Value = "163.87"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("D39").Value = "26.80"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "6.77"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("E42").Value = "  +5.37%  "
$ws.Range("D43").Value = "25.86"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").Value = "41.29"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "2.639.09"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "340.37"
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("D49").Value = "32.20"
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  -0.19%  "
